$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the current column C (NCTId), shifting
# NCTId..intervention_type (old C..L) one position to the right (new D..M).
$ws.Columns.Item(3).Insert()

# Header for the newly inserted column C.
$ws.Range("C1").Value2 = "statut_name"

# Map each row's existing statut_label (column B) to a human readable
# statut_name that is written into the new column C.
$map = @{
    "noir"   = "pas de résultat ni de publication"
    "rouge"  = "résultat et / ou publication posté"
    "orange" = "résultat et / ou publication posté dans les 36 mois"
    "vert"   = "résultat et / ou publication posté dans les 12 mois"
}

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $label = $ws.Cells.Item($r, 2).Value2
    if ($map.ContainsKey($label)) {
        $statutName = $map[$label]
    } else {
        $statutName = ""
    }
    $ws.Cells.Item($r, 3).Value2 = $statutName
}
